# New PO forecast model
# Updates three sheets: "Weekly Quantity", "Monthly Trend", "PO Forecast"
# with refreshed forecast data (new weeks/months appended, PO Forecast
# series recomputed going forward).

$wb = $excel.ActiveWorkbook
$dateFmt = "YYYY-MM-DD HH:MM:SS"

# ---------------------------------------------------------------------
# Sheet "Weekly Quantity": append 3 new weekly rows (A4:B6)
# ---------------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$weeklyNewRows = @(
    @(45669.99999999999, 5),
    @(45676.99999999999, 16),
    @(45683.99999999999, 7)
)
$r = 4
foreach ($row in $weeklyNewRows) {
    $wsWeekly.Cells.Item($r, 1).Value = $row[0]
    $wsWeekly.Cells.Item($r, 1).NumberFormat = $dateFmt
    $wsWeekly.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Sheet "Monthly Trend": append 1 new monthly row (A3:B3)
# ---------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Cells.Item(3, 1).Value = 45688.99999999999
$wsMonthly.Cells.Item(3, 1).NumberFormat = $dateFmt
$wsMonthly.Cells.Item(3, 2).Value = 28

# ---------------------------------------------------------------------
# Sheet "PO Forecast": rebuild the forecast series (A2:B14)
# ---------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("PO Forecast")
$forecastRows = @(
    @(45613.99999999999, 76),
    @(45620.99999999999, 69),
    @(45669.99999999999, 16),
    @(45676.99999999999, 8),
    @(45683.99999999999, 1),
    @(45690.99999999999, 0),
    @(45697.99999999999, 0),
    @(45704.99999999999, 0),
    @(45711.99999999999, 0),
    @(45718.99999999999, 0),
    @(45725.99999999999, 0),
    @(45732.99999999999, 0),
    @(45739.99999999999, 0)
)
$r = 2
foreach ($row in $forecastRows) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 1).NumberFormat = $dateFmt
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}
